$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 241.2
$ws.Range("J55").Value = 205.58333
$ws.Range("L55").Value = 205.58333
$ws.Range("N55").Value = -633.5833299999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("N69").Value = 0
$ws.Range("L69").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("N72").Value = 0
$ws.Range("L72").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2826.875
$ws.Range("I106").Value = 2806.1428
$ws.Range("K106").Value = 2806.1428
$ws.Range("M106").Value = -2175.1428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 75253.5
$ws.Range("J123").Value = 75253.5
$ws.Range("L123").Value = 75253.5
$ws.Range("N123").Value = -85053.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1658.2554
$ws.Range("I132").Value = 1569.0222
$ws.Range("K132").Value = 4707.0666
$ws.Range("M132").Value = -2177.0666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").Value = 0
$ws.Range("L133").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4882.5713
$ws.Range("I137").Value = 4985.375
$ws.Range("K137").Value = 14956.125
$ws.Range("M137").Value = -12406.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 65000
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2122776.2
$ws.Range("I32").Value = 2721968.5
$ws.Range("K32").Value = 2721968.5
$ws.Range("M32").Value = -2721681.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3014
$ws.Range("I45").Value = 2936
$ws.Range("K45").Value = 2936
$ws.Range("M45").Value = -2559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5360.364
$ws.Range("I61").Value = 2666.8838
$ws.Range("K61").Value = 2666.8838
$ws.Range("M61").Value = -2454.8838

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 58883.5
$ws.Range("J112").Value = 58883.5
$ws.Range("L112").Value = 58883.5
$ws.Range("N112").Value = -61837.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 10846.5
$ws.Range("I122").Value = 16852.6
$ws.Range("K122").Value = 50557.8
$ws.Range("M122").Value = -48107.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3521353
$ws.Range("J132").Value = 9804
$ws.Range("L132").Value = 29412
$ws.Range("N132").Value = -34472

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5360.364
$ws.Range("I136").Value = 2666.8838
$ws.Range("K136").Value = 8000.651400000001
$ws.Range("M136").Value = -5450.651400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8774319
$ws.Range("I20").Value = 13890346
$ws.Range("J20").Value = 3986.4285
$ws.Range("K20").Value = 13890346
$ws.Range("L20").Value = 3986.4285
$ws.Range("M20").Value = -13890099
$ws.Range("N20").Value = -4480.4285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4127.375
$ws.Range("I105").Value = 1881.6
$ws.Range("K105").Value = 1881.6
$ws.Range("M105").Value = -134.5999999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H131").Value = 40000
$ws.Range("J131").Value = 40000
$ws.Range("L131").Value = 40000
$ws.Range("N131").Value = -50080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 148.7
$ws.Range("I7").Value = 100.73333
$ws.Range("J7").Value = 292.6
$ws.Range("K7").Value = 100.73333
$ws.Range("L7").Value = 292.6
$ws.Range("M7").Value = 12.26667
$ws.Range("N7").Value = -518.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5634.1777
$ws.Range("I31").Value = 2468.68
$ws.Range("J31").Value = 9591.049999999999
$ws.Range("K31").Value = 2468.68
$ws.Range("L31").Value = 9591.049999999999
$ws.Range("M31").Value = -2173.68
$ws.Range("N31").Value = -10181.05

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5634.1777
$ws.Range("I34").Value = 2468.68
$ws.Range("J34").Value = 9591.049999999999
$ws.Range("K34").Value = 2468.68
$ws.Range("L34").Value = 9591.049999999999
$ws.Range("M34").Value = -2266.68
$ws.Range("N34").Value = -9995.049999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9339.764999999999
$ws.Range("I58").Value = 2551.6667
$ws.Range("J58").Value = 13042.363
$ws.Range("K58").Value = 2551.6667
$ws.Range("L58").Value = 13042.363
$ws.Range("M58").Value = -2348.6667
$ws.Range("N58").Value = -13448.363

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 8636.182000000001
$ws.Range("I132").Value = 3167.1667
$ws.Range("K132").Value = 9501.500100000001
$ws.Range("M132").Value = -6971.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 9339.764999999999
$ws.Range("I136").Value = 2551.6667
$ws.Range("J136").Value = 13042.363
$ws.Range("K136").Value = 7655.000100000001
$ws.Range("L136").Value = 39127.089
$ws.Range("M136").Value = -5105.000100000001
$ws.Range("N136").Value = -44227.089

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6813.3335
$ws.Range("J34").Value = 10172.625
$ws.Range("L34").Value = 30517.875
$ws.Range("N34").Value = -30685.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 48151732
$ws.Range("J55").Value = 33342250
$ws.Range("L55").Value = 100026750
$ws.Range("N55").Value = -100027104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3845.842
$ws.Range("I68").Value = 2183.6667
$ws.Range("J68").Value = 4613
$ws.Range("K68").Value = 6551.000100000001
$ws.Range("L68").Value = 13839
$ws.Range("M68").Value = -5740.000100000001
$ws.Range("N68").Value = -15461

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 3845.842
$ws.Range("I71").Value = 2183.6667
$ws.Range("J71").Value = 4613
$ws.Range("K71").Value = 19653.0003
$ws.Range("L71").Value = 41517
$ws.Range("M71").Value = -15597.0003
$ws.Range("N71").Value = -49629

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 2989.6667
$ws.Range("J74").Value = 2989.6667
$ws.Range("L74").Value = 8969.000100000001
$ws.Range("N74").Value = -11091.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H77").Value = 2989.6667
$ws.Range("J77").Value = 2989.6667
$ws.Range("L77").Value = 26907.0003
$ws.Range("N77").Value = -37515.0003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2278.913
$ws.Range("I107").Value = 623
$ws.Range("J107").Value = 2863.353
$ws.Range("K107").Value = 1869
$ws.Range("L107").Value = 8590.059000000001
$ws.Range("M107").Value = 51
$ws.Range("N107").Value = -12430.059

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2087.5715
$ws.Range("I113").Value = 1550
$ws.Range("K113").Value = 1550
$ws.Range("M113").Value = 620

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1918.3684
$ws.Range("I16").Value = 1913.8334
$ws.Range("K16").Value = 1913.8334
$ws.Range("M16").Value = -1743.8334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5167.3335
$ws.Range("I100").Value = 4250
$ws.Range("J100").Value = 7002
$ws.Range("K100").Value = 4250
$ws.Range("L100").Value = 7002
$ws.Range("M100").Value = -3709
$ws.Range("N100").Value = -8084

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5889.517
$ws.Range("I122").Value = 4639.55
$ws.Range("J122").Value = 8667.223
$ws.Range("K122").Value = 13918.65
$ws.Range("L122").Value = 26001.669
$ws.Range("M122").Value = -11468.65
$ws.Range("N122").Value = -30901.669

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7206.7407
$ws.Range("J132").Value = 14089.2
$ws.Range("L132").Value = 42267.60000000001
$ws.Range("N132").Value = -47327.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 56975.5
$ws.Range("J99").Value = 56975.5
$ws.Range("L99").Value = 56975.5
$ws.Range("N99").Value = -62965.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 38196
$ws.Range("J101").Value = 38196
$ws.Range("L101").Value = 38196
$ws.Range("N101").Value = -44686

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 28679.824
$ws.Range("I136").Value = 1834.5385
$ws.Range("J136").Value = 78535.36
$ws.Range("K136").Value = 5503.6155
$ws.Range("L136").Value = 235606.08
$ws.Range("M136").Value = -2953.6155
$ws.Range("N136").Value = -240706.08
